$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: period headers (D8:M8) ---
$ws.Range("D8").Value = "9 ماهه منتهی به 1399/09"
$ws.Range("E8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("F8").Value = "3 ماهه منتهی به 1400/03"
$ws.Range("G8").Value = "6 ماهه منتهی به 1400/06"
$ws.Range("H8").Value = "9 ماهه منتهی به 1400/09"
$ws.Range("I8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("J8").Value = "3 ماهه منتهی به 1401/03"
$ws.Range("K8").Value = "6 ماهه منتهی به 1401/06"
$ws.Range("L8").Value = "9 ماهه منتهی به 1401/09"
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: publish dates (D9:M9) ---
$ws.Range("D9").Value = "1401-01-15 (3)"
$ws.Range("E9").Value = "1401-03-11 (10)"
$ws.Range("F9").Value = "1401-04-30 (3)"
$ws.Range("G9").Value = "1401-08-18 (4)"
$ws.Range("H9").Value = "1401-10-29 (3)"
$ws.Range("I9").Value = "1402-02-27 (7)"
$ws.Range("J9").Value = "1401-04-30"
$ws.Range("K9").Value = "1401-08-18 (2)"
$ws.Range("L9").Value = "1401-10-29"
$ws.Range("M9").Value = "1402-02-27"

# --- Row 11 ---
$ws.Range("D11").Value = 1946051
$ws.Range("E11").Value = 2440996
$ws.Range("F11").Value = 957603
$ws.Range("G11").Value = 2187050
$ws.Range("H11").Value = 3489442
$ws.Range("I11").Value = 4297311
$ws.Range("J11").Value = 1210413
$ws.Range("K11").Value = 2660078
$ws.Range("L11").Value = 6355195
$ws.Range("M11").Value = 9819805

# --- Row 12 ---
$ws.Range("D12").Value = -1031407
$ws.Range("E12").Value = -1200366
$ws.Range("F12").Value = -466387
$ws.Range("G12").Value = -1028528
$ws.Range("H12").Value = -1756643
$ws.Range("I12").Value = -2127640
$ws.Range("J12").Value = -602665
$ws.Range("K12").Value = -1233793
$ws.Range("L12").Value = -3589415
$ws.Range("M12").Value = -5906404

# --- Row 13 ---
$ws.Range("D13").Value = 914644
$ws.Range("E13").Value = 1240630
$ws.Range("F13").Value = 491216
$ws.Range("G13").Value = 1158522
$ws.Range("H13").Value = 1732799
$ws.Range("I13").Value = 2169671
$ws.Range("J13").Value = 607748
$ws.Range("K13").Value = 1426285
$ws.Range("L13").Value = 2765780
$ws.Range("M13").Value = 3913401

# --- Row 14 ---
$ws.Range("D14").Value = -62807
$ws.Range("E14").Value = -116636
$ws.Range("F14").Value = -38454
$ws.Range("G14").Value = -82628
$ws.Range("H14").Value = -117786
$ws.Range("I14").Value = -177907
$ws.Range("J14").Value = -104526
$ws.Range("K14").Value = -109925
$ws.Range("L14").Value = -171308
$ws.Range("M14").Value = -213869

# --- Row 16 ---
$ws.Range("D16").Value = 16493
$ws.Range("E16").Value = 15630
$ws.Range("F16").Value = 1171
$ws.Range("G16").Value = 18186
$ws.Range("H16").Value = 29464
$ws.Range("I16").Value = 40427
$ws.Range("J16").Value = 9
$ws.Range("K16").Value = 5667
$ws.Range("L16").Value = 23767
$ws.Range("M16").Value = 28534

# --- Row 17 ---
$ws.Range("D17").Value = 868330
$ws.Range("E17").Value = 1139624
$ws.Range("F17").Value = 453933
$ws.Range("G17").Value = 1094080
$ws.Range("H17").Value = 1644477
$ws.Range("I17").Value = 2032191
$ws.Range("J17").Value = 503231
$ws.Range("K17").Value = 1322027
$ws.Range("L17").Value = 2618239
$ws.Range("M17").Value = 3728066

# --- Row 18 ---
$ws.Range("D18").Value = -88877
$ws.Range("E18").Value = -155571
$ws.Range("F18").Value = -27405
$ws.Range("G18").Value = -95101
$ws.Range("H18").Value = -160240
$ws.Range("I18").Value = -232392
$ws.Range("J18").Value = -25989
$ws.Range("K18").Value = -162382
$ws.Range("L18").Value = -233958
$ws.Range("M18").Value = -337504

# --- Row 19 ---
$ws.Range("D19").Value = 514
$ws.Range("E19").Value = -69982
$ws.Range("F19").Value = 23
$ws.Range("G19").Value = -30914
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = -33612
$ws.Range("J19").Value = 8510
$ws.Range("K19").Value = 8133
$ws.Range("L19").Value = 8846
$ws.Range("M19").Value = -2354

# --- Row 20 ---
$ws.Range("D20").Value = 779967
$ws.Range("E20").Value = 914071
$ws.Range("F20").Value = 426551
$ws.Range("G20").Value = 968065
$ws.Range("H20").Value = 1484237
$ws.Range("I20").Value = 1766187
$ws.Range("J20").Value = 485752
$ws.Range("K20").Value = 1167778
$ws.Range("L20").Value = 2393127
$ws.Range("M20").Value = 3388208

# --- Row 21 ---
$ws.Range("D21").Value = -177190
$ws.Range("E21").Value = -110014
$ws.Range("F21").Value = -96966
$ws.Range("G21").Value = -216904
$ws.Range("H21").Value = -309142
$ws.Range("I21").Value = -248862
$ws.Range("J21").Value = -109294
$ws.Range("K21").Value = -262729
$ws.Range("L21").Value = -538440
$ws.Range("M21").Value = -501836

# --- Row 22 ---
$ws.Range("D22").Value = 602777
$ws.Range("E22").Value = 804057
$ws.Range("F22").Value = 329585
$ws.Range("G22").Value = 751161
$ws.Range("H22").Value = 1175095
$ws.Range("I22").Value = 1517325
$ws.Range("J22").Value = 376458
$ws.Range("K22").Value = 905049
$ws.Range("L22").Value = 1854687
$ws.Range("M22").Value = 2886372

# --- Row 23 ---
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = 0

# --- Row 24 ---
$ws.Range("D24").Value = 602777
$ws.Range("E24").Value = 804057
$ws.Range("F24").Value = 329585
$ws.Range("G24").Value = 751161
$ws.Range("H24").Value = 1175095
$ws.Range("I24").Value = 1517325
$ws.Range("J24").Value = 376458
$ws.Range("K24").Value = 905049
$ws.Range("L24").Value = 1854687
$ws.Range("M24").Value = 2886372

# --- Row 25 ---
$ws.Range("D25").Value = 2364
$ws.Range("E25").Value = 3153
$ws.Range("F25").Value = 724
$ws.Range("G25").Value = 1651
$ws.Range("H25").Value = 2583
$ws.Range("I25").Value = 3335
$ws.Range("J25").Value = 827
$ws.Range("K25").Value = 1989
$ws.Range("L25").Value = 4076
$ws.Range("M25").Value = 3207

# --- Row 26 ---
$ws.Range("D26").Value = 255000
$ws.Range("E26").Value = 255000
$ws.Range("F26").Value = 455000
$ws.Range("G26").Value = 455000
$ws.Range("H26").Value = 455000
$ws.Range("I26").Value = 455000
$ws.Range("J26").Value = 455000
$ws.Range("K26").Value = 455000
$ws.Range("L26").Value = 455000
$ws.Range("M26").Value = 900000

# --- Row 27 ---
$ws.Range("D27").Value = 670
$ws.Range("E27").Value = 893
$ws.Range("F27").Value = 366
$ws.Range("G27").Value = 835
$ws.Range("H27").Value = 1306
$ws.Range("I27").Value = 1686
$ws.Range("J27").Value = 418
$ws.Range("K27").Value = 1006
$ws.Range("L27").Value = 2061
$ws.Range("M27").Value = 3207

# --- Column widths: shift pattern left by one (D..M) ---
$ws.Columns.Item(4).ColumnWidth = 27.166666666666668
$ws.Columns.Item(5).ColumnWidth = 28.166666666666668
$ws.Columns.Item(6).ColumnWidth = 27.166666666666668
$ws.Columns.Item(7).ColumnWidth = 27.166666666666668
$ws.Columns.Item(8).ColumnWidth = 27.166666666666668
$ws.Columns.Item(9).ColumnWidth = 28.166666666666668
$ws.Columns.Item(10).ColumnWidth = 27.166666666666668
$ws.Columns.Item(11).ColumnWidth = 27.166666666666668
$ws.Columns.Item(12).ColumnWidth = 27.166666666666668
$ws.Columns.Item(13).ColumnWidth = 28.166666666666668
